$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 1666678.6
$ws.Range("I6").Value = 1666678.6
$ws.Range("K6").Value = 5000035.800000001
$ws.Range("M6").Value = -4999923.800000001

$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("N40").ClearContents()

$ws.Range("H41").Value = 477.7143
$ws.Range("I41").Value = 340.25
$ws.Range("K41").Value = 340.25
$ws.Range("M41").Value = 99.75

$ws.Range("H52").Value = 25999
$ws.Range("I52").Value = 1498
$ws.Range("J52").Value = 38249.5
$ws.Range("K52").Value = 4494
$ws.Range("L52").Value = 114748.5
$ws.Range("M52").Value = -4334
$ws.Range("N52").Value = -115068.5

$ws.Range("H92").Value = 319
$ws.Range("I92").Value = 315
$ws.Range("J92").Value = 325
$ws.Range("K92").Value = 315
$ws.Range("L92").Value = 325
$ws.Range("M92").Value = 933
$ws.Range("N92").Value = -2821

$ws.Range("H99").Value = 2671.25
$ws.Range("I99").Value = 2671.25
$ws.Range("K99").Value = 8013.75
$ws.Range("M99").Value = -6515.75

$ws.Range("H107").Value = 693
$ws.Range("I107").Value = 413.75
$ws.Range("K107").Value = 413.75
$ws.Range("M107").Value = 1506.25

$ws.Range("H123").Value = 109296.664
$ws.Range("J123").Value = 109296.664
$ws.Range("L123").Value = 109296.664
$ws.Range("N123").Value = -119096.664

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 14428
$ws.Range("I45").Value = 8749
$ws.Range("K45").Value = 8749
$ws.Range("M45").Value = -8372

$ws.Range("H74").Value = 3459
$ws.Range("I74").Value = 2498.3333
$ws.Range("K74").Value = 2498.3333
$ws.Range("M74").Value = -1624.3333

$ws.Range("H77").Value = 3459
$ws.Range("I77").Value = 2498.3333
$ws.Range("K77").Value = 12491.6665
$ws.Range("M77").Value = -8123.666499999999

$ws.Range("H110").Value = 1336.3334
$ws.Range("I110").Value = 1336.3334
$ws.Range("K110").Value = 1336.3334
$ws.Range("M110").Value = 708.6666

$ws.Range("H113").Value = 55555
$ws.Range("J113").Value = 55555
$ws.Range("L113").Value = 55555
$ws.Range("N113").Value = -64233

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 4495
$ws.Range("I99").Value = 4495
$ws.Range("K99").Value = 4495
$ws.Range("M99").Value = -2997

$ws.Range("H107").Value = 4399.6
$ws.Range("I107").Value = 4332.6665
$ws.Range("J107").Value = 4500
$ws.Range("K107").Value = 4332.6665
$ws.Range("L107").Value = 4500
$ws.Range("M107").Value = -2412.6665
$ws.Range("N107").Value = -8340

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6324.5386
$ws.Range("I31").Value = 1798
$ws.Range("J31").Value = 7682.5
$ws.Range("K31").Value = 1798
$ws.Range("L31").Value = 7682.5
$ws.Range("M31").Value = -1503
$ws.Range("N31").Value = -8272.5

$ws.Range("H34").Value = 6324.5386
$ws.Range("I34").Value = 1798
$ws.Range("J34").Value = 7682.5
$ws.Range("K34").Value = 1798
$ws.Range("L34").Value = 7682.5
$ws.Range("M34").Value = -1596
$ws.Range("N34").Value = -8086.5

$ws.Range("H92").Value = 45000
$ws.Range("J92").Value = 45000
$ws.Range("L92").Value = 45000
$ws.Range("N92").Value = -49992

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 41204668
$ws.Range("I4").Value = 290.2
$ws.Range("J4").Value = 68674260
$ws.Range("K4").Value = 870.5999999999999
$ws.Range("L4").Value = 206022780
$ws.Range("M4").Value = -758.5999999999999
$ws.Range("N4").Value = -206023004

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()

$ws.Range("H48").Value = 0
$ws.Range("I48").Value = 0
$ws.Range("K48").Value = 0
$ws.Range("M48").ClearContents()

$ws.Range("H70").Value = 6400
$ws.Range("J70").Value = 6400
$ws.Range("L70").Value = 6400
$ws.Range("N70").Value = -6940

$ws.Range("H73").Value = 6400
$ws.Range("J73").Value = 6400
$ws.Range("L73").Value = 6400
$ws.Range("N73").Value = -8272

$ws.Range("H97").Value = 1099.5
$ws.Range("I97").Value = 1099.5
$ws.Range("K97").Value = 1099.5
$ws.Range("M97").Value = -603.5

$ws.Range("H102").Value = 3374.6667
$ws.Range("I102").Value = 3374.6667
$ws.Range("K102").Value = 3374.6667
$ws.Range("M102").Value = -1752.6667

$ws.Range("H113").Value = 4766.5
$ws.Range("I113").Value = 5266.3335
$ws.Range("K113").Value = 5266.3335
$ws.Range("M113").Value = -3096.3335

$ws.Range("H132").Value = 8241.733
$ws.Range("I132").Value = 6125.077
$ws.Range("K132").Value = 18375.231
$ws.Range("M132").Value = -15845.231

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H26").Value = 2000
$ws.Range("I26").Value = 2000
$ws.Range("K26").Value = 2000
$ws.Range("M26").Value = -1705

$ws.Range("H46").Value = 999
$ws.Range("I46").Value = 999
$ws.Range("K46").Value = 999
$ws.Range("M46").Value = -811

$ws.Range("H55").Value = 500
$ws.Range("I55").Value = 500
$ws.Range("K55").Value = 500
$ws.Range("M55").Value = -327

$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("M68").ClearContents()

$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("M71").ClearContents()

$ws.Range("H132").Value = 14719.75
$ws.Range("I132").Value = 5002
$ws.Range("J132").Value = 17959
$ws.Range("K132").Value = 15006
$ws.Range("L132").Value = 53877
$ws.Range("M132").Value = -12476
$ws.Range("N132").Value = -58937

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1033
$ws.Range("I96").Value = 900
$ws.Range("J96").Value = 1099.5
$ws.Range("K96").Value = 900
$ws.Range("L96").Value = 1099.5
$ws.Range("M96").Value = 473
$ws.Range("N96").Value = -3845.5

$ws.Range("H105").Value = 40233
$ws.Range("J105").Value = 40233
$ws.Range("L105").Value = 40233
$ws.Range("N105").Value = -47221

$ws.Range("H113").Value = 929.3333
$ws.Range("I113").Value = 929.3333
$ws.Range("K113").Value = 2787.9999
$ws.Range("M113").Value = -617.9998999999998

$ws.Range("H126").Value = 2000
$ws.Range("I126").Value = 2000
$ws.Range("K126").Value = 6000
$ws.Range("M126").Value = -3530

$ws.Range("H132").Value = 8683.416999999999
$ws.Range("I132").Value = 7743
$ws.Range("K132").Value = 23229
$ws.Range("M132").Value = -20699
